# Appends the "Aprilis 27." progress-log entry (Leaderboard view section)
# at the very end of the document body, right before the final section break,
# by inserting a literal WordprocessingML fragment via Range.InsertXML.
$d = $word.ActiveDocument

$endPos = $d.Content.End
$insertionPoint = $d.Range($endPos, $endPos)

$xmlFragment = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Április </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>27.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Custom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> színek hozzáadása az </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>asset-ekhez</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, osztály létrehozása, amelyen keresztül elérhetőek a különböző </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>View</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> objektumok </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>modifier-eiből</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Ezeket a színeket felhasználva esztétikusabb megjelenést alakítottam ki az alkalmazásban megjelenő gomboknak, köztük a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Form</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> objektumokon belül szereplő gomboknak is.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Ranglista nézet</w:t></w:r><w:r><w:t xml:space="preserve"> („</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LeaderboardView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”)</w:t></w:r><w:r><w:t xml:space="preserve"> kialakítása: az oldal tetején egy </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Picker</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> található </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>„.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>segmented</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">” stílusban, ezzel választható ki, hogy a bérlők/foglalók vagy a kiadók ranglistájára kíváncsi a felhasználó. Amikor kiválasztotta, akkor egyelőre egy példa </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>usereket</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tartalmazó tömb alapján kilistázza az alkalmazás a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>userek</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> adatait képpel, névvel és az adott ranglistához tartozó pontszámmal.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>userek</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ranglista adataihoz külön </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>View</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> objektumot</w:t></w:r><w:r><w:t xml:space="preserve"> („</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>UserLeaderboardView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”)</w:t></w:r><w:r><w:t xml:space="preserve"> hoztam létre. Ez megjeleníti a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> profilképét, amely </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>default</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> esetben egy </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>system</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> kép, a felhasználó nevét és pontszámát. </w:t></w:r><w:r><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>View</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> egy </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Hstack-ből</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> áll, amelyben baloldalon található a profilkép, és jobboldalon egy </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Vstack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, amelyben egymás alatt találhatóak a nevet és a pontszámot tartalmazó Text objektumok.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>userek</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> adatait egy külön erre a célra elkészített objektum alapján szerzi meg a ranglistán szereplő </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> adatokat megjelenítő </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>View</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Ebben szerepelnek a korábban már felsorolt adatok, illetve egy UUID, amely alapján azonosítani lehet az objektumokat (és így pl. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ForEach</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> objektumban lehet őket használni).</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Elvégzendő feladat:</w:t></w:r><w:r><w:t xml:space="preserve"> a példa </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>usereket</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tartalmazó tömböket kiszervezni egy külön objektumba, akár egy </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>environment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> változóba, és a pontszám alapján történő </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sorbarendezést</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> elvégezni rajtuk. Erre a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>View</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> megjelenítésekor, illetve módosításakor nincs lehetőség, így külön objektum kell majd erre a célra. Illetve, a későbbiekben, amikor nem példa tömbökkel fog működni az alkalmazás, akkor is könnyebb külön karbantartani az erre a célra létrehozott objektumot és adott esetben frissíteni azt. (Tervezési minták alkalmazása.)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xmlFragment)
